$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The output in column C had a gap (…8, 9, 11, 12…) because two rows were
# missing from the data set. Re-insert them so C becomes continuous 1..15.

# Missing row: "Oranges" belongs right after the "Pears" row (was row 4).
$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = 42102.37480324074
$ws.Range("B4").Value = "Oranges"
$ws.Range("C4").Value = 4

# Missing row: "S " belongs right after the "B " row of the first cycle
# (now row 10, after the previous insert shifted everything down by one).
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "S "
$ws.Range("C10").Value = 10
